$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.140.83"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.762.38"
$ws.Range("E3").Value = "  -2.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.22"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9973"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3784"
$ws.Range("E7").Value = "  -3.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3367"
$ws.Range("E8").Value = "  -3.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.83"
$ws.Range("E9").Value = "  -5.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.133"
$ws.Range("E10").Value = "  -5.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07217"
$ws.Range("E11").Value = "  -4.89%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.65"
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("B13").Value = "BinanceUSD"
$ws.Range("C13").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9976"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.208"
$ws.Range("E14").Value = "  -4.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.216"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "1.758.36"
$ws.Range("E16").Value = "  -2.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001058"
$ws.Range("E17").Value = "  -4.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06600"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.83"
$ws.Range("E19").Value = "  -4.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9978"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.05"
$ws.Range("E21").Value = "  -3.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.292"
$ws.Range("E22").Value = "  -4.35%  "
$ws.Range("D23").Value = "28.106.60"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.73"
$ws.Range("E24").Value = "  -5.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.393"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.51"
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.92"
$ws.Range("E27").Value = "  -6.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.351"
$ws.Range("E28").Value = "  -7.06%  "
$ws.Range("D29").Value = "1.958.96"
$ws.Range("E29").Value = "  -2.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.271"
$ws.Range("E30").Value = "  -15.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "131.95"
$ws.Range("E31").Value = "  -3.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.010"
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.840"
$ws.Range("E33").Value = "  -5.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08814"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.33"
$ws.Range("E35").Value = "  -5.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02349"
$ws.Range("E36").Value = "  -3.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6645"
$ws.Range("E37").Value = "  -4.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.190"
$ws.Range("E38").Value = "  -5.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06211"
$ws.Range("E39").Value = "  -5.28%  "
$ws.Range("E40").Value = "  -4.33%  "
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.458"
$ws.Range("E42").Value = "  -9.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.060"
$ws.Range("E43").Value = "  -5.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9969"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.78"
$ws.Range("E45").Value = "  -5.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6079"
$ws.Range("E46").Value = "  -5.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.815"
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.90"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.023"
$ws.Range("E49").Value = "  -6.06%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07216"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.186"
$ws.Range("E51").Value = "  +2.59%  "
